$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid and Absent become 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Absent becomes 1
$ws.Range("H4").Value = 1

# Row 5: Absent becomes 1
$ws.Range("H5").Value = 1

# Row 6: Total Attendance Count and Real become 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-18: Absent becomes 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
